$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grade values (5) for the specified cells
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 5

$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5

$ws.Range("F21").Value = 5

$ws.Range("F26").Value = 5

$ws.Range("F27").Value = 5

$ws.Range("D29").Value = 5

$ws.Range("D32").Value = 5
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = 5

# Update the active selection to F32
$ws.Range("F32").Select()
